$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 41775.3
$ws.Range("J86").Value = 68900.664
$ws.Range("L86").Value = 68900.664
$ws.Range("N86").Value = -71146.664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 41775.3
$ws.Range("J89").Value = 68900.664
$ws.Range("L89").Value = 344503.32
$ws.Range("N89").Value = -355735.32

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 909.60785
$ws.Range("J129").Value = 914.7447
$ws.Range("L129").Value = 2744.2341
$ws.Range("N129").Value = -12744.2341

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1342.2
$ws.Range("I107").Value = 1070.3334
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 1070.3334
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 849.6666
$ws.Range("N107").Value = -5590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6839.425
$ws.Range("I31").Value = 1415.8485
$ws.Range("J31").Value = 32407.715
$ws.Range("K31").Value = 1415.8485
$ws.Range("L31").Value = 32407.715
$ws.Range("M31").Value = -1120.8485
$ws.Range("N31").Value = -32997.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6839.425
$ws.Range("I34").Value = 1415.8485
$ws.Range("J34").Value = 32407.715
$ws.Range("K34").Value = 1415.8485
$ws.Range("L34").Value = 32407.715
$ws.Range("M34").Value = -1213.8485
$ws.Range("N34").Value = -32811.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2512.96
$ws.Range("I58").Value = 2143.7144
$ws.Range("K58").Value = 2143.7144
$ws.Range("M58").Value = -1940.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 655.7692
$ws.Range("I107").Value = 625
$ws.Range("J107").Value = 705
$ws.Range("K107").Value = 625
$ws.Range("L107").Value = 705
$ws.Range("M107").Value = 1295
$ws.Range("N107").Value = -4545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3847.077
$ws.Range("I134").Value = 3453.0908
$ws.Range("K134").Value = 10359.2724
$ws.Range("M134").Value = -7824.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2512.96
$ws.Range("I136").Value = 2143.7144
$ws.Range("K136").Value = 6431.1432
$ws.Range("M136").Value = -3881.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 766.58826
$ws.Range("I34").Value = 105.5
$ws.Range("J34").Value = 970
$ws.Range("K34").Value = 316.5
$ws.Range("L34").Value = 2910
$ws.Range("M34").Value = -232.5
$ws.Range("N34").Value = -3078

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 195.57143
$ws.Range("I40").Value = 117.25
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 469
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -400
$ws.Range("N40").Value = -1338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2963.3333
$ws.Range("J48").Value = 2963.3333
$ws.Range("L48").Value = 8889.999899999999
$ws.Range("N48").Value = -9389.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4701.857
$ws.Range("I63").Value = 3231
$ws.Range("J63").Value = 6663
$ws.Range("K63").Value = 9693
$ws.Range("L63").Value = 19989
$ws.Range("M63").Value = -8944
$ws.Range("N63").Value = -21487

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2070.6667
$ws.Range("I64").Value = 212
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 636
$ws.Range("L64").Value = 9000
$ws.Range("M64").Value = -366
$ws.Range("N64").Value = -9540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 4701.857
$ws.Range("I66").Value = 3231
$ws.Range("J66").Value = 6663
$ws.Range("K66").Value = 29079
$ws.Range("L66").Value = 59967
$ws.Range("M66").Value = -25335
$ws.Range("N66").Value = -67455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2070.6667
$ws.Range("I67").Value = 212
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 636
$ws.Range("L67").Value = 9000
$ws.Range("M67").Value = 300
$ws.Range("N67").Value = -10872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1362.1666
$ws.Range("I70").Value = 634.6
$ws.Range("K70").Value = 1903.8
$ws.Range("M70").Value = -1588.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 1362.1666
$ws.Range("I73").Value = 634.6
$ws.Range("K73").Value = 1903.8
$ws.Range("M73").Value = -811.8000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5323.3335
$ws.Range("I80").Value = 4500
$ws.Range("J80").Value = 5488
$ws.Range("K80").Value = 13500
$ws.Range("L80").Value = 16464
$ws.Range("M80").Value = -12564
$ws.Range("N80").Value = -18336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 103143.8
$ws.Range("I82").Value = 613
$ws.Range("J82").Value = 128776.5
$ws.Range("K82").Value = 1839
$ws.Range("L82").Value = 386329.5
$ws.Range("M82").Value = -1433
$ws.Range("N82").Value = -387141.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5323.3335
$ws.Range("I83").Value = 4500
$ws.Range("J83").Value = 5488
$ws.Range("K83").Value = 40500
$ws.Range("L83").Value = 49392
$ws.Range("M83").Value = -35820
$ws.Range("N83").Value = -58752

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 103143.8
$ws.Range("I85").Value = 613
$ws.Range("J85").Value = 128776.5
$ws.Range("K85").Value = 1839
$ws.Range("L85").Value = 386329.5
$ws.Range("M85").Value = -435
$ws.Range("N85").Value = -389137.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 15529.538
$ws.Range("I87").Value = 7032
$ws.Range("J87").Value = 19306.223
$ws.Range("K87").Value = 21096
$ws.Range("L87").Value = 57918.66900000001
$ws.Range("M87").Value = -19848
$ws.Range("N87").Value = -60414.66900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 15529.538
$ws.Range("I90").Value = 7032
$ws.Range("J90").Value = 19306.223
$ws.Range("K90").Value = 63288
$ws.Range("L90").Value = 173756.007
$ws.Range("M90").Value = -57048
$ws.Range("N90").Value = -186236.007

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 539.0909
$ws.Range("I114").Value = 364.5
$ws.Range("J114").Value = 638.8570999999999
$ws.Range("K114").Value = 1093.5
$ws.Range("L114").Value = 1916.5713
$ws.Range("M114").Value = 2160.5
$ws.Range("N114").Value = -8424.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1712.0476
$ws.Range("I129").Value = 1282.5
$ws.Range("J129").Value = 1813.1177
$ws.Range("K129").Value = 3847.5
$ws.Range("L129").Value = 5439.3531
$ws.Range("M129").Value = 1152.5
$ws.Range("N129").Value = -15439.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1123.321
$ws.Range("J131").Value = 1227.5916
$ws.Range("L131").Value = 3682.7748
$ws.Range("N131").Value = -13762.7748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2157.6667
$ws.Range("I136").Value = 1989.875
$ws.Range("K136").Value = 5969.625
$ws.Range("M136").Value = -869.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4751.4
$ws.Range("I80").Value = 6266.4
$ws.Range("J80").Value = 3236.4
$ws.Range("K80").Value = 6266.4
$ws.Range("L80").Value = 3236.4
$ws.Range("M80").Value = -5268.4
$ws.Range("N80").Value = -5232.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4751.4
$ws.Range("I83").Value = 6266.4
$ws.Range("J83").Value = 3236.4
$ws.Range("K83").Value = 31332
$ws.Range("L83").Value = 16182
$ws.Range("M83").Value = -26340
$ws.Range("N83").Value = -26166

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1890.6666
$ws.Range("I107").Value = 2673.375
$ws.Range("J107").Value = 325.25
$ws.Range("K107").Value = 2673.375
$ws.Range("L107").Value = 325.25
$ws.Range("M107").Value = -753.375
$ws.Range("N107").Value = -4165.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H129").Value = 42422.848
$ws.Range("J129").Value = 42422.848
$ws.Range("L129").Value = 42422.848
$ws.Range("N129").Value = -52422.848

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 14851.667
$ws.Range("J3").Value = 14851.667
$ws.Range("L3").Value = 14851.667
$ws.Range("N3").Value = -15075.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 14851.667
$ws.Range("J15").Value = 14851.667
$ws.Range("L15").Value = 14851.667
$ws.Range("N15").Value = -15191.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 42110.77
$ws.Range("I40").Value = 52959
$ws.Range("J40").Value = 5950
$ws.Range("K40").Value = 52959
$ws.Range("L40").Value = 5950
$ws.Range("M40").Value = -52823
$ws.Range("N40").Value = -6222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 51218.3
$ws.Range("I46").Value = 143714
$ws.Range("J46").Value = 1412.9231
$ws.Range("K46").Value = 143714
$ws.Range("L46").Value = 1412.9231
$ws.Range("M46").Value = -143526
$ws.Range("N46").Value = -1788.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 2425
$ws.Range("I107").Value = 2425
$ws.Range("K107").Value = 2425
$ws.Range("M107").Value = -505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = $null
